$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row updates: column letter -> new value, per row number
$updates = @{
  2 = @{ 'D' = '24.941.98'; 'E' = '  +2.23%  ' }
  3 = @{ 'D' = '1.678.28'; 'E' = '  +1.84%  ' }
  4 = @{ 'D' = '1.002'; 'E' = '  -0.13%  ' }
  5 = @{ 'D' = '328.91'; 'E' = '  +7.11%  ' }
  6 = @{ 'D' = '1.000'; 'E' = '  -0.04%  ' }
  7 = @{ 'D' = '0.3662'; 'E' = '  +1.37%  ' }
  8 = @{ 'D' = '46.99'; 'E' = '  -1.37%  ' }
  9 = @{ 'D' = '0.3253'; 'E' = '  -0.52%  ' }
  10 = @{ 'E' = '  +2.59%  ' }
  11 = @{ 'D' = '0.07090'; 'E' = '  +2.53%  ' }
  12 = @{ 'D' = '1.000'; 'E' = '  +0.06%  ' }
  13 = @{ 'D' = '6.100'; 'E' = '  +3.08%  ' }
  14 = @{ 'D' = '19.69'; 'E' = '  +2.68%  ' }
  15 = @{ 'D' = '1.680.42'; 'E' = '  +2.14%  ' }
  16 = @{ 'D' = '6.644'; 'E' = '  +1.15%  ' }
  17 = @{ 'D' = '0.00001051'; 'E' = '  +1.21%  ' }
  18 = @{ 'D' = '0.06590'; 'E' = '  +1.45%  ' }
  19 = @{ 'D' = '0.9998'; 'E' = '  +0.02%  ' }
  20 = @{ 'D' = '78.95'; 'E' = '  +3.50%  ' }
  21 = @{ 'D' = '15.94'; 'E' = '  +1.98%  ' }
  22 = @{ 'D' = '5.927'; 'E' = '  +0.40%  ' }
  23 = @{ 'D' = '12.91'; 'E' = '  +6.15%  ' }
  24 = @{ 'D' = '24.943.49'; 'E' = '  +2.43%  ' }
  25 = @{ 'D' = '2.449'; 'E' = '  +1.10%  ' }
  26 = @{ 'D' = '2.418'; 'E' = '  +4.22%  ' }
  27 = @{ 'D' = '148.03'; 'E' = '  +1.37%  ' }
  28 = @{ 'D' = '18.75'; 'E' = '  +3.04%  ' }
  29 = @{ 'D' = '1.863.94'; 'E' = '  +1.86%  ' }
  30 = @{ 'D' = '125.84'; 'E' = '  +1.32%  ' }
  31 = @{ 'D' = '1.185'; 'E' = '  +1.88%  ' }
  32 = @{ 'D' = '4.075'; 'E' = '  +0.82%  ' }
  33 = @{ 'D' = '5.776'; 'E' = '  +3.46%  ' }
  34 = @{ 'D' = '0.08502'; 'E' = '  +2.13%  ' }
  35 = @{ 'D' = '1.644'; 'E' = '  -1.69%  ' }
  36 = @{ 'E' = '  +0.77%  ' }
  37 = @{ 'D' = '5.193'; 'E' = '  -0.14%  ' }
  38 = @{ 'D' = '0.02255'; 'E' = '  +2.80%  ' }
  39 = @{ 'B' = 'Algorand'; 'C' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'; 'D' = '0.2098'; 'E' = '  +2.84%  ' }
  40 = @{ 'B' = 'TrustWalletToken'; 'C' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; 'D' = '1.228'; 'E' = '  +2.01%  ' }
  41 = @{ 'B' = 'Hedera'; 'C' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; 'D' = '0.06020'; 'E' = '  +0.27%  ' }
  42 = @{ 'D' = '8.254'; 'E' = '  +1.29%  ' }
  43 = @{ 'D' = '0.9993'; 'E' = '  -0.06%  ' }
  44 = @{ 'D' = '0.5964'; 'E' = '  +3.11%  ' }
  45 = @{ 'D' = '13.77'; 'E' = '  +9.62%  ' }
  46 = @{ 'D' = '3.844'; 'E' = '  +3.24%  ' }
  47 = @{ 'E' = '  +3.93%  ' }
  48 = @{ 'E' = '  +3.33%  ' }
  49 = @{ 'D' = '1.968'; 'E' = '  +2.05%  ' }
  50 = @{ 'D' = '0.07026'; 'E' = '  +2.03%  ' }
  51 = @{ 'D' = '1.192'; 'E' = '  +3.61%  ' }
}

# Cells whose new value would be mis-detected as a number by COM
# (e.g. "1.002", "328.91") must be pre-formatted as Text so they are
# written back as strings, matching the source inlineStr cells.
$textForceCells = @(
  'D4', 'D5', 'D6', 'D7', 'D8', 'D9', 'D11', 'D12', 'D13', 'D14',
  'D16', 'D17', 'D18', 'D19', 'D20', 'D21', 'D22', 'D23', 'D25', 'D26',
  'D27', 'D28', 'D30', 'D31', 'D32', 'D33', 'D34', 'D35', 'D37', 'D38',
  'D39', 'D40', 'D41', 'D42', 'D43', 'D44', 'D45', 'D46', 'D49', 'D50',
  'D51'
)

foreach ($addr in $textForceCells) {
  $ws.Range($addr).NumberFormat = "@"
}

foreach ($rowNum in $updates.Keys) {
  $rowData = $updates[$rowNum]
  foreach ($col in $rowData.Keys) {
    $ws.Range("$col$rowNum").Value = $rowData[$col]
  }
}
